$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (41 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 849.72
$ws.Range("I15").Value = 849.72
$ws.Range("K15").Value = 2549.16
$ws.Range("M15").Value = -2380.16
$ws.Range("H51").Value = 10000
$ws.Range("I51").Value = 10000
$ws.Range("K51").Value = 10000
$ws.Range("M51").Value = -9516
$ws.Range("H53").Value = 2245.2104
$ws.Range("I53").Value = 247
$ws.Range("K53").Value = 247
$ws.Range("M53").Value = 390
$ws.Range("H112").Value = 4445515.5
$ws.Range("J112").Value = 4445515.5
$ws.Range("L112").Value = 13336546.5
$ws.Range("N112").Value = -13338762.5
$ws.Range("H129").Value = 239230.14
$ws.Range("J129").Value = 264382
$ws.Range("L129").Value = 793146
$ws.Range("N129").Value = -803146
$ws.Range("H132").Value = 2514.175
$ws.Range("I132").Value = 2584.5
$ws.Range("J132").Value = 2115.6667
$ws.Range("K132").Value = 7753.5
$ws.Range("L132").Value = 6347.000100000001
$ws.Range("M132").Value = -5223.5
$ws.Range("N132").Value = -11407.0001
$ws.Range("H137").Value = 1719.7567
$ws.Range("I137").Value = 1564.3667
$ws.Range("J137").Value = 2385.7144
$ws.Range("K137").Value = 4693.1001
$ws.Range("L137").Value = 7157.1432
$ws.Range("M137").Value = -2143.1001
$ws.Range("N137").Value = -12257.1432
$ws.Range("H138").Value = 10991163
$ws.Range("I138").Value = 21277424
$ws.Range("J138").Value = 3565.9092
$ws.Range("K138").Value = 63832272
$ws.Range("L138").Value = 10697.7276
$ws.Range("M138").Value = -63827132
$ws.Range("N138").Value = -20977.7276

# --- Sheet: ARM (23 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1050.58
$ws.Range("I32").Value = 1066.9348
$ws.Range("J32").Value = 862.5
$ws.Range("K32").Value = 1066.9348
$ws.Range("L32").Value = 862.5
$ws.Range("M32").Value = -779.9348
$ws.Range("N32").Value = -1436.5
$ws.Range("H61").Value = 429717.5
$ws.Range("I61").Value = 530591.9
$ws.Range("K61").Value = 530591.9
$ws.Range("M61").Value = -530379.9
$ws.Range("H74").Value = 2359.1292
$ws.Range("I74").Value = 2271.6667
$ws.Range("K74").Value = 2271.6667
$ws.Range("M74").Value = -1397.6667
$ws.Range("H77").Value = 2359.1292
$ws.Range("I77").Value = 2271.6667
$ws.Range("K77").Value = 11358.3335
$ws.Range("M77").Value = -6990.333500000001
$ws.Range("H136").Value = 429717.5
$ws.Range("I136").Value = 530591.9
$ws.Range("K136").Value = 1591775.7
$ws.Range("M136").Value = -1589225.7

# --- Sheet: BSM (19 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1550.1904
$ws.Range("I86").Value = 1379.2084
$ws.Range("K86").Value = 1379.2084
$ws.Range("M86").Value = -256.2084
$ws.Range("H89").Value = 1550.1904
$ws.Range("I89").Value = 1379.2084
$ws.Range("K89").Value = 6896.041999999999
$ws.Range("M89").Value = -1280.041999999999
$ws.Range("H99").Value = 1757.2727
$ws.Range("I99").Value = 1907.5
$ws.Range("K99").Value = 1907.5
$ws.Range("M99").Value = -409.5
$ws.Range("H134").Value = 2862.54
$ws.Range("I134").Value = 3000.3171
$ws.Range("J134").Value = 2234.889
$ws.Range("K134").Value = 9000.951300000001
$ws.Range("L134").Value = 6704.667
$ws.Range("M134").Value = -6465.951300000001
$ws.Range("N134").Value = -11774.667

# --- Sheet: CRP (37 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3055
$ws.Range("I31").Value = 1758.742
$ws.Range("K31").Value = 1758.742
$ws.Range("M31").Value = -1463.742
$ws.Range("H34").Value = 3055
$ws.Range("I34").Value = 1758.742
$ws.Range("K34").Value = 1758.742
$ws.Range("M34").Value = -1556.742
$ws.Range("H58").Value = 15777.529
$ws.Range("I58").Value = 1038.2069
$ws.Range("J58").Value = 101265.6
$ws.Range("K58").Value = 1038.2069
$ws.Range("L58").Value = 101265.6
$ws.Range("M58").Value = -835.2068999999999
$ws.Range("N58").Value = -101671.6
$ws.Range("H99").Value = 20837022
$ws.Range("J99").Value = 35718052
$ws.Range("L99").Value = 35718052
$ws.Range("N99").Value = -35721048
$ws.Range("H126").Value = 20837022
$ws.Range("J126").Value = 35718052
$ws.Range("L126").Value = 107154156
$ws.Range("N126").Value = -107159096
$ws.Range("H134").Value = 939.9535
$ws.Range("I134").Value = 856.6177
$ws.Range("J134").Value = 1254.7778
$ws.Range("K134").Value = 2569.8531
$ws.Range("L134").Value = 3764.3334
$ws.Range("M134").Value = -34.85310000000027
$ws.Range("N134").Value = -8834.3334
$ws.Range("H136").Value = 15777.529
$ws.Range("I136").Value = 1038.2069
$ws.Range("J136").Value = 101265.6
$ws.Range("K136").Value = 3114.620699999999
$ws.Range("L136").Value = 303796.8
$ws.Range("M136").Value = -564.6206999999995
$ws.Range("N136").Value = -308896.8

# --- Sheet: CUL (36 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1765062
$ws.Range("I4").Value = 325.35715
$ws.Range("J4").Value = 10000500
$ws.Range("K4").Value = 976.0714499999999
$ws.Range("L4").Value = 30001500
$ws.Range("M4").Value = -864.0714499999999
$ws.Range("N4").Value = -30001724
$ws.Range("H23").Value = 405.6
$ws.Range("I23").Value = 55.25
$ws.Range("J23").Value = 639.1667
$ws.Range("K23").Value = 165.75
$ws.Range("L23").Value = 1917.5001
$ws.Range("M23").Value = 69.25
$ws.Range("N23").Value = -2387.5001
$ws.Range("H109").Value = 3618.45
$ws.Range("I109").Value = 1125.2222
$ws.Range("J109").Value = 5658.364
$ws.Range("K109").Value = 3375.6666
$ws.Range("L109").Value = 16975.092
$ws.Range("M109").Value = -2335.6666
$ws.Range("N109").Value = -19055.092
$ws.Range("H129").Value = 209412.45
$ws.Range("I129").Value = 852.0769
$ws.Range("J129").Value = 455892.9
$ws.Range("K129").Value = 2556.2307
$ws.Range("L129").Value = 1367678.7
$ws.Range("M129").Value = 2443.7693
$ws.Range("N129").Value = -1377678.7
$ws.Range("H131").Value = 721.90814
$ws.Range("J131").Value = 731.7957
$ws.Range("L131").Value = 2195.3871
$ws.Range("N131").Value = -12275.3871
$ws.Range("H139").Value = 1590.625
$ws.Range("J139").Value = 3964.1667
$ws.Range("L139").Value = 11892.5001
$ws.Range("N139").Value = -22172.5001

# --- Sheet: GSM (21 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3729.35
$ws.Range("I80").Value = 2954.4285
$ws.Range("J80").Value = 4146.615
$ws.Range("K80").Value = 2954.4285
$ws.Range("L80").Value = 4146.615
$ws.Range("M80").Value = -1956.4285
$ws.Range("N80").Value = -6142.615
$ws.Range("H83").Value = 3729.35
$ws.Range("I83").Value = 2954.4285
$ws.Range("J83").Value = 4146.615
$ws.Range("K83").Value = 14772.1425
$ws.Range("L83").Value = 20733.075
$ws.Range("M83").Value = -9780.1425
$ws.Range("N83").Value = -30717.075
$ws.Range("H132").Value = 19095.6
$ws.Range("I132").Value = 2797.88
$ws.Range("J132").Value = 100584.2
$ws.Range("K132").Value = 8393.639999999999
$ws.Range("L132").Value = 301752.6
$ws.Range("M132").Value = -5863.639999999999
$ws.Range("N132").Value = -306812.6

# --- Sheet: LTW (25 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1093750
$ws.Range("I2").Value = 1178571.4
$ws.Range("J2").Value = 499999.5
$ws.Range("K2").Value = 1178571.4
$ws.Range("L2").Value = 499999.5
$ws.Range("M2").Value = -1178459.4
$ws.Range("N2").Value = -500223.5
$ws.Range("H82").Value = 2270.1333
$ws.Range("I82").Value = 2032
$ws.Range("J82").Value = 2925
$ws.Range("K82").Value = 2032
$ws.Range("L82").Value = 2925
$ws.Range("M82").Value = -1671
$ws.Range("N82").Value = -3647
$ws.Range("H85").Value = 2270.1333
$ws.Range("I85").Value = 2032
$ws.Range("J85").Value = 2925
$ws.Range("K85").Value = 2032
$ws.Range("L85").Value = 2925
$ws.Range("M85").Value = -784
$ws.Range("N85").Value = -5421
$ws.Range("H136").Value = 1386.3334
$ws.Range("I136").Value = 1279.1052
$ws.Range("K136").Value = 3837.3156
$ws.Range("M136").Value = -1287.3156

# --- Sheet: WVR (29 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 38001.5
$ws.Range("J2").Value = 38001.5
$ws.Range("L2").Value = 38001.5
$ws.Range("N2").Value = -38225.5
$ws.Range("H81").Value = 52632670
$ws.Range("I81").Value = 1098.9333
$ws.Range("J81").Value = 250001070
$ws.Range("K81").Value = 2197.8666
$ws.Range("L81").Value = 500002140
$ws.Range("M81").Value = -1136.8666
$ws.Range("N81").Value = -500004262
$ws.Range("H84").Value = 52632670
$ws.Range("I84").Value = 1098.9333
$ws.Range("J84").Value = 250001070
$ws.Range("K84").Value = 10989.333
$ws.Range("L84").Value = 2500010700
$ws.Range("M84").Value = -5685.332999999999
$ws.Range("N84").Value = -2500021308
$ws.Range("H132").Value = 765.8214
$ws.Range("I132").Value = 765.8214
$ws.Range("K132").Value = 2297.4642
$ws.Range("M132").Value = 232.5357999999997
$ws.Range("H136").Value = 16668695
$ws.Range("I136").Value = 24391220
$ws.Range("J136").Value = 4297.8423
$ws.Range("K136").Value = 73173660
$ws.Range("L136").Value = 12893.5269
$ws.Range("M136").Value = -73171110
$ws.Range("N136").Value = -17993.5269
